# Generate Report for Handback
#
# Fills in the "Latest Target File" / "Latest Handback File" /
# "Latest Handback DateTime" columns for the zh-cn and de-de handback
# tables, flips the Overview "Status" columns from "Ready for handoff"
# to "Handed back: in sync with en-US", and widens a few columns that
# now hold longer content.

$wb = $excel.ActiveWorkbook

$baseUrl = "https://github.com/OpenLocalizationTestOrg/ol-test1/blob/363282c018352b6f48372b14ee248381f3a3a76a/e2e/f22cedc0-676f-45f3-9b46-105fd66ba8c8.md"
$mdName  = "f22cedc0-676f-45f3-9b46-105fd66ba8c8.md"

# ---------------------------------------------------------------
# Overview sheet: status for both locales flips to "handed back"
# ---------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "Handed back: in sync with en-US"
$overview.Range("F2").Value = "Handed back: in sync with en-US"

# Columns E/F need to be wider to fit the new, longer status text.
$overview.Columns.Item(5).ColumnWidth = 29.166666666666668
$overview.Columns.Item(6).ColumnWidth = 29.166666666666668

# ---------------------------------------------------------------
# zh-cn handback table
# ---------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

# Status column (C) flips to "handed back" too, same as the Overview
# sheet, and is now longer -> widen it.
$zhcn.Range("C2").Value = "Handed back: in sync with en-US"
$zhcn.Columns.Item(3).ColumnWidth = 29.166666666666668

# Latest Target File (J) / Latest Handback File (K) / Latest Handback
# DateTime (L) are now populated with the generated handback report.
$zhcn.Columns.Item(10).ColumnWidth = 39.166666666666664
$zhcn.Columns.Item(11).ColumnWidth = 39.166666666666664

$zhcn.Hyperlinks.Add($zhcn.Range("J2"), $baseUrl, "", "", $mdName)
# Match the look of the workbook's other hyperlinks (underlined, cornflowerblue).
$zhcn.Range("J2").Font.Underline = $true
$zhcn.Range("J2").Font.Color = 15570276
$zhcn.Range("K2").Value = "f22cedc0-676f-45f3-9b46-105fd66ba8c8.486070234ff1031d98e40249a22f43ffc20b7c62.zh-cn.xlf"
$zhcn.Range("L2").Value = "2017-01-03 04:20:36"

# ---------------------------------------------------------------
# de-de handback table
# ---------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("C2").Value = "Handed back: in sync with en-US"
$dede.Columns.Item(3).ColumnWidth = 29.166666666666668
$dede.Columns.Item(10).ColumnWidth = 39.166666666666664
$dede.Columns.Item(11).ColumnWidth = 39.166666666666664

$dede.Hyperlinks.Add($dede.Range("J2"), $baseUrl, "", "", $mdName)
$dede.Range("J2").Font.Underline = $true
$dede.Range("J2").Font.Color = 15570276
$dede.Range("K2").Value = "f22cedc0-676f-45f3-9b46-105fd66ba8c8.486070234ff1031d98e40249a22f43ffc20b7c62.de-de.xlf"
$dede.Range("L2").Value = "2017-01-03 04:20:47"
